$d = $word.ActiveDocument

# --- Locate the target paragraph ("Baz changes") precisely, so the
#     replacements below don't also match the unrelated "changes" in
#     an earlier paragraph. ---
$p = $d.Paragraphs.Item(5)

# Set paragraph spacing to single (w:line="240" w:lineRule="auto")
$p.Range.ParagraphFormat.LineSpacingRule = 0

# Replace "Baz chan" -> greeting text, scoped to this paragraph only.
$r = $p.Range
$r.Find.Execute("Baz chan", $true, $false, $false, $false, $false, $true, 1, $false, "Hi, here is the changed made by dnguye22. My git hub account is minhduc", 2)

# Re-fetch the paragraph (text length changed) and replace "ges" -> student id text.
$p = $d.Paragraphs.Item(5)
$r2 = $p.Range
$r2.Find.Execute("ges", $true, $false, $false, $false, $false, $true, 1, $false, "512, student id:22240998", 2)

# Remove one of the two trailing empty paragraphs.
$empty = $d.Paragraphs.Item(6)
$empty.Range.Delete()
